$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.5019577983990334, 0.5019577983990334, 0.3105127224087574, 0.5572366843709748, 0.251107588723997, 14),
    @(0.3307979378781269, 0.3343319746521877, 0.1645692088094062, 0.4056713063668741, 0.244411730890862, 13),
    @(0.2377396291553046, 0.2452017974494533, 0.09171268429751143, 0.3028410214906683, 0.1959384495933066, 12),
    @(0.373625158811909, 0.373625158811909, 0.1693451806315346, 0.4115157112815191, 0.1808987657993869, 11),
    @(0.3392289017490621, 0.3392289017490621, 0.1435050013861421, 0.3788205398155466, 0.1777287371382058, 10),
    @(0.323658031409521, 0.323658031409521, 0.1411500930911699, 0.3756994717738766, 0.202348754060147, 9),
    @(0.3534922319734252, 0.3534922319734252, 0.1628866010887078, 0.4035921221836568, 0.2082027666157719, 8),
    @(0.3545941574938994, 0.3545941574938994, 0.1673332494109491, 0.4090638695985617, 0.2202928468255327, 7),
    @(0.3114747142229594, 0.3114747142229594, 0.1227287858089793, 0.3503266844089661, 0.1756551902177846, 6),
    @(0.3600284238150658, 0.3600284238150658, 0.1554279835465356, 0.3942435586620733, 0.1796090114379531, 5)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
    $ws.Cells.Item($row, 6).Value = $data[$i][4]
    $ws.Cells.Item($row, 7).Value = $data[$i][5]
}
